$wb = $excel.ActiveWorkbook

# ALC!row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 117
$ws.Range("I12").Value = 118.8
$ws.Range("J12").Value = 112.5
$ws.Range("K12").Value = 118.8
$ws.Range("L12").Value = 112.5
$ws.Range("M12").Value = 51.2
$ws.Range("N12").Value = -452.5

# ALC!row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9998.666999999999
$ws.Range("J43").Value = 9999.666999999999
$ws.Range("L43").Value = 9999.666999999999
$ws.Range("N43").Value = -10137.667

# ALC!row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 891.625
$ws.Range("I55").Value = 840.8333
$ws.Range("K55").Value = 840.8333
$ws.Range("M55").Value = -626.8333

# ALC!row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4099.1665
$ws.Range("I64").Value = 4099.1665
$ws.Range("K64").Value = 4099.1665
$ws.Range("M64").Value = -3851.1665

# ALC!row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4099.1665
$ws.Range("I67").Value = 4099.1665
$ws.Range("K67").Value = 4099.1665
$ws.Range("M67").Value = -3241.1665

# ALC!row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5250.6
$ws.Range("J88").Value = 6967.6665
$ws.Range("L88").Value = 6967.6665
$ws.Range("N88").Value = -7779.6665

# ALC!row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 5250.6
$ws.Range("J91").Value = 6967.6665
$ws.Range("L91").Value = 6967.6665
$ws.Range("N91").Value = -9775.666499999999

# ALC!row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 499
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("M103").Value = 1497
$ws.Range("N103").Value = -2669

# ARM!row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1474
$ws.Range("I2").Value = 1474
$ws.Range("K2").Value = 1474
$ws.Range("M2").Value = -1361

# ARM!row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2412.6667
$ws.Range("I88").Value = 2419
$ws.Range("K88").Value = 2419
$ws.Range("M88").Value = -2013

# ARM!row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2412.6667
$ws.Range("I91").Value = 2419
$ws.Range("K91").Value = 2419
$ws.Range("M91").Value = -1015

# ARM!row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3236.25
$ws.Range("I97").Value = 3270
$ws.Range("K97").Value = 3270
$ws.Range("M97").Value = -2774

# ARM!row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1474
$ws.Range("I116").Value = 1474
$ws.Range("K116").Value = 1474
$ws.Range("M116").Value = 820

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1466.3334
$ws.Range("I122").Value = 1337.125
$ws.Range("K122").Value = 4011.375
$ws.Range("M122").Value = -1561.375

# ARM!row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6974.9165
$ws.Range("I132").Value = 4911
$ws.Range("K132").Value = 14733
$ws.Range("M132").Value = -12203

# BSM!row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1474
$ws.Range("I3").Value = 1474
$ws.Range("K3").Value = 1474
$ws.Range("M3").Value = -1360

# CRP!row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# CRP!row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 29225
$ws.Range("J88").Value = 29225
$ws.Range("L88").Value = 29225
$ws.Range("N88").Value = -30037

# CRP!row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# CRP!row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 29225
$ws.Range("J91").Value = 29225
$ws.Range("L91").Value = 29225
$ws.Range("N91").Value = -32033

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3714.2856
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# CRP!row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2466.3333
$ws.Range("I105").Value = 2466.3333
$ws.Range("K105").Value = 2466.3333
$ws.Range("M105").Value = -719.3332999999998

# CRP!row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3714.2856
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

# CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4735.8887
$ws.Range("I132").Value = 2374.7144
$ws.Range("K132").Value = 7124.1432
$ws.Range("M132").Value = -4594.1432

# CUL!row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 261.42856
$ws.Range("I29").Value = 182.75
$ws.Range("J29").Value = 366.33334
$ws.Range("K29").Value = 548.25
$ws.Range("L29").Value = 1099.00002
$ws.Range("M29").Value = -271.25
$ws.Range("N29").Value = -1653.00002

# CUL!row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752

# CUL!row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 999.6667
$ws.Range("I122").Value = 999.6667
$ws.Range("K122").Value = 2999.0001
$ws.Range("M122").Value = -549.0001000000002

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7200.3335
$ws.Range("I132").Value = 3202.2
$ws.Range("J132").Value = 12198
$ws.Range("K132").Value = 9606.599999999999
$ws.Range("L132").Value = 36594
$ws.Range("M132").Value = -7076.599999999999
$ws.Range("N132").Value = -41654

# LTW!row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -3888

# LTW!row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9530

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12685.333
$ws.Range("I132").Value = 6778
$ws.Range("K132").Value = 20334
$ws.Range("M132").Value = -17804

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9388.777
$ws.Range("I136").Value = 8999.5
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 26998.5
$ws.Range("L136").Value = 28500
$ws.Range("M136").Value = -24448.5
$ws.Range("N136").Value = -33600

# WVR!row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 875
$ws.Range("I107").Value = 770.8333
$ws.Range("K107").Value = 2312.4999
$ws.Range("M107").Value = -392.4998999999998

# WVR!row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2229
$ws.Range("I113").Value = 2229
$ws.Range("K113").Value = 6687
$ws.Range("M113").Value = -4517

# WVR!row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7565.8667
$ws.Range("I132").Value = 5408.091
$ws.Range("K132").Value = 16224.273
$ws.Range("M132").Value = -13694.273
